# Duplicate the data rows (row 3) on both sheets twice more (rows 4-6),
# incrementing the application-number column (C) each time, then update
# the active-sheet/selection state to match the post-edit workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # PC1_CitizenshipCertificate
$ws2 = $wb.Worksheets.Item(2)   # PC7_PassportRenewal

# --- Sheet1 (PC1_CitizenshipCertificate): clone row 3 into rows 4,5,6 ---
# Row 3 has no B3 cell, so copy A and C:V separately to avoid materialising
# an empty B cell.
for ($i = 4; $i -le 6; $i++) {
    $ws1.Range("A3").Copy($ws1.Range("A$i"))
    $ws1.Range("C3:V3").Copy($ws1.Range("C$i"))
    $ws1.Cells.Item($i, 3).Value = $i - 1
}

# --- Sheet2 (PC7_PassportRenewal): clone row 3 into rows 4,5,6 ---
# Row 3 has no D3 cell, so copy A:C and E:N separately.
for ($i = 4; $i -le 6; $i++) {
    $ws2.Range("A3:C3").Copy($ws2.Range("A$i"))
    $ws2.Range("E3:N3").Copy($ws2.Range("E$i"))
    $ws2.Cells.Item($i, 3).Value = $i - 1
}

# --- View state: sheet1 loses the tab selection / old cell selection,
#     sheet2 becomes the active tab with C2:C6 selected ---
$null = $ws1.Range("D12").Select()
$null = $ws2.Range("C2:C6").Select()
$null = $ws2.Activate()
